$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'325.45"
$ws.Range("E2").Value = "'-0.65%"
$ws.Range("D3").Value = "'44.28"
$ws.Range("E3").Value = "'-1.85%"
$ws.Range("D4").Value = "'5.496"
$ws.Range("E4").Value = "'-1.88%"
$ws.Range("D5").Value = "'0.07986"
$ws.Range("E5").Value = "'-1.81%"
$ws.Range("D6").Value = "'2.009"
$ws.Range("E6").Value = "'5.62%"
$ws.Range("D7").Value = "'4.300"
$ws.Range("E7").Value = "'-1.00%"
$ws.Range("B8").Value = "BTSEToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D8").Value = "'2.565"
$ws.Range("E8").Value = "'-6.87%"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "'0.9470"
$ws.Range("E9").Value = "'-0.78%"
$ws.Range("D10").Value = "'0.1144"
$ws.Range("E10").Value = "'-3.20%"
$ws.Range("D11").Value = "'0.1838"
$ws.Range("E11").Value = "'-3.89%"
$ws.Range("D12").Value = "'12.58"
$ws.Range("E12").Value = "'45.94%"
$ws.Range("D13").Value = "'0.09581"
$ws.Range("E13").Value = "'-3.30%"
$ws.Range("D14").Value = "'0.04559"
$ws.Range("E14").Value = "'8.76%"
$ws.Range("D15").Value = "'0.1064"
$ws.Range("E15").Value = "'-0.07%"
$ws.Range("D16").Value = "'0.001270"
$ws.Range("E16").Value = "'-0.99%"
$ws.Range("D17").Value = "'0.04074"
$ws.Range("E17").Value = "'-4.77%"
$ws.Range("D18").Value = "'0.005890"
$ws.Range("E18").Value = "'-0.47%"
$ws.Range("E19").Value = "'-6.30%"
$ws.Range("D20").Value = "'0.3479"
$ws.Range("E20").Value = "'-0.21%"
$ws.Range("D21").Value = "'0.1409"
$ws.Range("E21").Value = "'2.39%"
$ws.Range("D22").Value = "'0.2545"
$ws.Range("E22").Value = "'-2.09%"
$ws.Range("D23").Value = "'0.001240"
$ws.Range("E23").Value = "'-0.10%"
$ws.Range("D24").Value = "'0.004295"
$ws.Range("E24").Value = "'-6.80%"
$ws.Range("E25").Value = "'-3.73%"
$ws.Range("D26").Value = "'0.0003743"
$ws.Range("E26").Value = "'-6.55%"
$ws.Range("D38").Value = "'0.02532"
$ws.Range("E38").Value = "'-5.98%"
$ws.Range("D39").Value = "'0.05519"
$ws.Range("E39").Value = "'-2.31%"
$ws.Range("D40").Value = "'0.007533"
$ws.Range("E40").Value = "'-2.22%"
$ws.Range("D41").Value = "'0.1389"
$ws.Range("E41").Value = "'-0.80%"
$ws.Range("D42").Value = "'0.007587"
$ws.Range("E42").Value = "'-33.06%"
$ws.Range("D43").Value = "'0.002015"
$ws.Range("E43").Value = "'-2.36%"
$ws.Range("D44").Value = "'0.008383"
$ws.Range("E44").Value = "'-3.53%"
$ws.Range("D45").Value = "'0.00007100"
$ws.Range("E45").Value = "'-0.05%"
$ws.Range("E46").Value = "'-0.43%"
$ws.Range("E47").Value = "'0.98%"
$ws.Range("D48").Value = "'0.004223"
$ws.Range("E48").Value = "'21.50%"
$ws.Range("D49").Value = "'0.00002100"
$ws.Range("E49").Value = "'-0.43%"
$ws.Range("D50").Value = "'0.0002000"
$ws.Range("E50").Value = "'-0.43%"
